# Update the player/team assignments on the roster sheet.
# Only columns A (player) and C (team) change for the affected rows;
# the position in column B stays the same for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Stephen Curry"
$ws.Range("C2").Value = "Golden State Warriors"

$ws.Range("A7").Value = "Karl-Anthony Towns"
$ws.Range("C7").Value = "New York Knicks"

$ws.Range("A10").Value = "Jalen Duren"
$ws.Range("C10").Value = "Detroit Pistons"

$ws.Range("A11").Value = "Mark Williams"
$ws.Range("C11").Value = "Charlotte Hornets"

$ws.Range("A14").Value = "Tyrese Maxey"
$ws.Range("C14").Value = "Philadelphia 76ers"

$ws.Range("A15").Value = "Jarrett Allen"
$ws.Range("C15").Value = "Cleveland Cavaliers"

$ws.Range("A16").Value = "Daniel Gafford"
$ws.Range("C16").Value = "Dallas Mavericks"
